# DC-Colos.xlsx: insert a new colo row ("NNG" / Nanning, China) right before
# the existing "IAD" row (row 275), pushing every subsequent row down by one.
# The new row has no lat/lon data yet (matches the source diff, which leaves
# G/H blank for this colo).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 275 - everything currently at/after row 275
# (IAD, ATL, BOS, ... YHZ) shifts down by one row.
$ws.Rows.Item(275).Insert()

# Fill in the new row's data.
$ws.Range("A275").Value = "NNG"
$ws.Range("B275").Value = "Nanning, China"
$ws.Range("C275").Value = "Asia"
$ws.Range("D275").Value = "Nanning"
$ws.Range("E275").Value = "China"
$ws.Range("F275").Value = "CN"
# G275 / H275 (lat / lon) intentionally left blank - no coordinates yet.

# Re-apply the standard "colo code" cell formatting (bold, bordered,
# center/top aligned) to A275, matching every other row in the table -
# the plain row insert above doesn't carry that formatting over cleanly.
$ws.Range("A276").Copy()
$ws.Range("A275").PasteSpecial(-4122)
$excel.CutCopyMode = 0
